# Applies the "Adding vehicle scenarios part 1" change:
#  - CREDENTIALS sheet: the existing username/password pair is pushed down to
#    become a second credential pair (rows 6-7, "username.2"/"password.2"),
#    and a brand new credential pair is written into rows 1-2.
#  - URL sheet: the qa1 URLs are repointed at the dev environment.

$wb = $excel.ActiveWorkbook

# ======================================================================
# CREDENTIALS sheet
# ======================================================================
$ws = $wb.Worksheets.Item("CREDENTIALS")

# Stash current formatting of B1 (hyperlink style) and B2 (plain style) in
# scratch cells so it can be re-applied later without the engine minting
# brand new style records every time a hyperlink is (re)created.
$ws.Range("B1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the existing username/password pair down into the new row 6/7 slot.
$ws.Range("A6").Value = "username.2"
$ws.Range("B6").Value = $ws.Range("B1").Value()
$ws.Range("A7").Value = "password.2"
$ws.Range("B7").Value = $ws.Range("B2").Value()

# Re-apply formatting to the moved cells to match the original row1/row2
# look (A6/B6 mirror A1/B1, A7/B7 mirror A2/B2).
$ws.Range("A1").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Recreate the hyperlink for B6 (same mailto target the original B1
# hyperlink used), then restore its hyperlink-style formatting.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:kanabecadmin@mobile22.com")
$ws.Range("Z1").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the new credential pair into row 1/row 2.
$ws.Range("B1").Value = "sukanabec@mailinator.com"
$ws.Range("B2").Value = "changeMe22@"
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:sukanabec@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:changeMe22@")

# Restore hyperlink-style formatting on B1/B2 (Hyperlinks.Add mints a
# fresh generic style; put back the workbook's existing "Hipervinculo" xf).
$ws.Range("Z1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Clean up scratch cells.
$ws.Range("Z1").Clear()
$ws.Range("Z2").Clear()

# Selection on this sheet (will not remain the active tab - URL is
# reselected below, last, so it stays the active sheet on save).
$null = $ws.Range("B15").Select()

# ======================================================================
# URL sheet
# ======================================================================
$ws2 = $wb.Worksheets.Item("URL")

# Stash the existing hyperlink-cell formatting (style shared by B2:B4).
$ws2.Range("B2").Copy()
$ws2.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Repoint the three URLs at the dev environment. Hyperlinks.Delete() clears
# every hyperlink on the sheet, so all three are rebuilt together,
# preserving the B2->rId1, B3->rId2, B4->rId3 ordering.
$ws2.Hyperlinks.Delete()
$ws2.Range("B2").Value = "https://driver.dev.mobile22.com"
$ws2.Range("B3").Value = "https://vehicle.dev.mobile22.com"
$ws2.Range("B4").Value = "https://callcenter.dev.mobile22.com"
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://driver.dev.mobile22.com/")
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://vehicle.dev.mobile22.com/")
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://callcenter.dev.mobile22.com/")

# Restore the original hyperlink-cell formatting across B2:B4.
$ws2.Range("Z1").Copy()
$ws2.Range("B2:B4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws2.Range("Z1").Clear()

# Selection - select this last so URL remains the active tab on save.
$null = $ws2.Range("C11").Select()
